# Update Leve profit figures across sheets (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 897.2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 897.2
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2691.6
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3027.6
$ws.Range("H28").Value = 1467.5555
$ws.Range("J28").Value = 339
$ws.Range("L28").Value = 339
$ws.Range("N28").Value = -1309
$ws.Range("H40").Value = 1771.68
$ws.Range("I40").Value = 1487.3125
$ws.Range("K40").Value = 1487.3125
$ws.Range("M40").Value = -1312.3125
$ws.Range("H100").Value = 1351.6154
$ws.Range("I100").Value = 940.6667
$ws.Range("J100").Value = 2276.25
$ws.Range("K100").Value = 940.6667
$ws.Range("L100").Value = 2276.25
$ws.Range("M100").Value = -399.6667
$ws.Range("N100").Value = -3358.25
$ws.Range("H138").Value = 5085.815
$ws.Range("I138").Value = 2836.1667
$ws.Range("J138").Value = 5728.5713
$ws.Range("K138").Value = 8508.500100000001
$ws.Range("L138").Value = 17185.7139
$ws.Range("M138").Value = -3368.500100000001
$ws.Range("N138").Value = -27465.7139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1719.2
$ws.Range("I2").Value = 1774
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 1774
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -1661
$ws.Range("N2").Value = -1726
$ws.Range("H10").Value = 403180.8
$ws.Range("J10").Value = 4005
$ws.Range("L10").Value = 4005
$ws.Range("N10").Value = -4345
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H116").Value = 1719.2
$ws.Range("I116").Value = 1774
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 1774
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 520
$ws.Range("N116").Value = -6088
$ws.Range("H122").Value = 1075.7778
$ws.Range("I122").Value = 883.1429000000001
$ws.Range("K122").Value = 2649.4287
$ws.Range("M122").Value = -199.4287000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1719.2
$ws.Range("I3").Value = 1774
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 1774
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -1660
$ws.Range("N3").Value = -1728
$ws.Range("H22").Value = 228.4
$ws.Range("I22").Value = 160.5
$ws.Range("K22").Value = 160.5
$ws.Range("M22").Value = 12.5
$ws.Range("H64").Value = 548.25
$ws.Range("I64").Value = 689.6667
$ws.Range("K64").Value = 689.6667
$ws.Range("M64").Value = -464.6667
$ws.Range("H67").Value = 548.25
$ws.Range("I67").Value = 689.6667
$ws.Range("K67").Value = 689.6667
$ws.Range("M67").Value = 90.33330000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3649.5
$ws.Range("I62").Value = 3149
$ws.Range("J62").Value = 4150
$ws.Range("K62").Value = 3149
$ws.Range("L62").Value = 4150
$ws.Range("M62").Value = -2525
$ws.Range("N62").Value = -5398
$ws.Range("H65").Value = 3649.5
$ws.Range("I65").Value = 3149
$ws.Range("J65").Value = 4150
$ws.Range("K65").Value = 15745
$ws.Range("L65").Value = 20750
$ws.Range("M65").Value = -12625
$ws.Range("N65").Value = -26990
$ws.Range("H99").Value = 2117.1667
$ws.Range("I99").Value = 1863.875
$ws.Range("J99").Value = 2623.75
$ws.Range("K99").Value = 1863.875
$ws.Range("L99").Value = 2623.75
$ws.Range("M99").Value = -365.875
$ws.Range("N99").Value = -5619.75
$ws.Range("H122").Value = 6181
$ws.Range("I122").Value = 6649.1
$ws.Range("K122").Value = 19947.3
$ws.Range("M122").Value = -17497.3
$ws.Range("H126").Value = 2117.1667
$ws.Range("I126").Value = 1863.875
$ws.Range("J126").Value = 2623.75
$ws.Range("K126").Value = 5591.625
$ws.Range("L126").Value = 7871.25
$ws.Range("M126").Value = -3121.625
$ws.Range("N126").Value = -12811.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 105.333336
$ws.Range("I17").Value = 110.5
$ws.Range("J17").Value = 95
$ws.Range("K17").Value = 331.5
$ws.Range("L17").Value = 285
$ws.Range("M17").Value = -162.5
$ws.Range("N17").Value = -623
$ws.Range("H35").Value = 100
$ws.Range("J35").Value = 100
$ws.Range("L35").Value = 300
$ws.Range("N35").Value = -876
$ws.Range("H39").Value = 5000
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H55").Value = 750
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354
$ws.Range("M55").Value = -1323
$ws.Range("H107").Value = 343
$ws.Range("I107").Value = 343
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1029
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 891
$ws.Range("N107").ClearContents()
$ws.Range("H110").Value = 926.6667
$ws.Range("I110").Value = 926.6667
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2780.0001
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1309.9999
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 1168.4667
$ws.Range("J113").Value = 1254.3334
$ws.Range("L113").Value = 3763.0002
$ws.Range("N113").Value = -8103.0002
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("H138").Value = 1819.8334
$ws.Range("I138").Value = 1819.8334
$ws.Range("K138").Value = 5459.5002
$ws.Range("M138").Value = -319.5002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H107").Value = 2860.7273
$ws.Range("I107").Value = 1295
$ws.Range("J107").Value = 4165.5
$ws.Range("K107").Value = 1295
$ws.Range("L107").Value = 4165.5
$ws.Range("M107").Value = 625
$ws.Range("N107").Value = -8005.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5739.6
$ws.Range("I40").Value = 5566
$ws.Range("K40").Value = 5566
$ws.Range("M40").Value = -5430
$ws.Range("J68").Value = 1500
$ws.Range("L68").Value = 1500
$ws.Range("N68").Value = -2998
$ws.Range("J71").Value = 1500
$ws.Range("L71").Value = 7500
$ws.Range("N71").Value = -14988
$ws.Range("H122").Value = 6555.222
$ws.Range("J122").Value = 8713.571
$ws.Range("L122").Value = 26140.713
$ws.Range("N122").Value = -31040.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3396.6667
$ws.Range("J6").Value = 3345.6667
$ws.Range("L6").Value = 3345.6667
$ws.Range("N6").Value = -3575.6667
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H96").Value = 3601.5
$ws.Range("I96").Value = 3003
$ws.Range("J96").Value = 4200
$ws.Range("K96").Value = 3003
$ws.Range("L96").Value = 4200
$ws.Range("N96").Value = -6946
$ws.Range("M96").Value = -1630
$ws.Range("H122").Value = 1883.3334
$ws.Range("I122").Value = 1883.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5650.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3200.0002
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3207.182
$ws.Range("I132").Value = 3207.182
$ws.Range("K132").Value = 9621.545999999998
$ws.Range("M132").Value = -7091.545999999998
$ws.Range("H140").Value = 89994.5
$ws.Range("J140").Value = 89994.5
$ws.Range("L140").Value = 89994.5
$ws.Range("N140").Value = -100354.5
$ws.Range("H141").Value = 231999.67
$ws.Range("J141").Value = 231999.67
$ws.Range("L141").Value = 231999.67
